# Auto-generated Excel COM-interop script applying the Kujata_Profits market-price update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3263.9092
$ws.Range("I74").Value = 3400.5
$ws.Range("K74").Value = 3400.5
$ws.Range("M74").Value = -2464.5

$ws.Range("H77").Value = 3263.9092
$ws.Range("I77").Value = 3400.5
$ws.Range("K77").Value = 17002.5
$ws.Range("M77").Value = -12322.5

$ws.Range("H88").Value = 1124988.4
$ws.Range("I88").Value = 996.5
$ws.Range("J88").Value = 1374764.4
$ws.Range("K88").Value = 996.5
$ws.Range("L88").Value = 1374764.4
$ws.Range("M88").Value = -590.5
$ws.Range("N88").Value = -1375576.4

$ws.Range("H91").Value = 1124988.4
$ws.Range("I91").Value = 996.5
$ws.Range("J91").Value = 1374764.4
$ws.Range("K91").Value = 996.5
$ws.Range("L91").Value = 1374764.4
$ws.Range("M91").Value = 407.5
$ws.Range("N91").Value = -1377572.4

$ws.Range("H92").Value = 1918
$ws.Range("I92").Value = 1652.6666
$ws.Range("K92").Value = 1652.6666
$ws.Range("M92").Value = -404.6666

$ws.Range("H96").Value = 2051.1
$ws.Range("I96").Value = 2126.375
$ws.Range("J96").Value = 1750
$ws.Range("K96").Value = 6379.125
$ws.Range("L96").Value = 5250
$ws.Range("M96").Value = -5006.125
$ws.Range("N96").Value = -7996

$ws.Range("H100").Value = 1427.7858
$ws.Range("I100").Value = 1367.5555
$ws.Range("J100").Value = 1536.2
$ws.Range("K100").Value = 1367.5555
$ws.Range("L100").Value = 1536.2
$ws.Range("M100").Value = -826.5554999999999
$ws.Range("N100").Value = -2618.2

$ws.Range("H112").Value = 3587.375
$ws.Range("J112").Value = 3942.7144
$ws.Range("L112").Value = 11828.1432
$ws.Range("N112").Value = -14044.1432

$ws.Range("H124").Value = 39933.332
$ws.Range("J124").Value = 39933.332
$ws.Range("L124").Value = 39933.332
$ws.Range("N124").Value = -49753.332

$ws.Range("H132").Value = 8038.294
$ws.Range("I132").Value = 6323.2
$ws.Range("K132").Value = 18969.6
$ws.Range("M132").Value = -16439.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H32").Value = 5719.298
$ws.Range("I32").Value = 5626.2607
$ws.Range("K32").Value = 5626.2607
$ws.Range("M32").Value = -5339.2607

$ws.Range("H132").Value = 2725.5
$ws.Range("I132").Value = 2377.2222
$ws.Range("J132").Value = 3770.3333
$ws.Range("K132").Value = 7131.6666
$ws.Range("L132").Value = 11310.9999
$ws.Range("M132").Value = -4601.6666
$ws.Range("N132").Value = -16370.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 337.5
$ws.Range("I22").Value = 283.33334
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 283.33334
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -110.33334
$ws.Range("N22").Value = -846

$ws.Range("H59").Value = 50000
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

$ws.Range("H94").Value = 27778456
$ws.Range("I94").Value = 35714870
$ws.Range("J94").Value = 1005
$ws.Range("K94").Value = 35714870
$ws.Range("L94").Value = 1005
$ws.Range("M94").Value = -35714419
$ws.Range("N94").Value = -1907

$ws.Range("H105").Value = 90910350
$ws.Range("I105").Value = 90910350
$ws.Range("K105").Value = 90910350
$ws.Range("M105").Value = -90908603

$ws.Range("H123").Value = 37949.5
$ws.Range("J123").Value = 40253.332
$ws.Range("L123").Value = 40253.332
$ws.Range("N123").Value = -50053.332

$ws.Range("H134").Value = 8851.571
$ws.Range("I134").Value = 1156.5454
$ws.Range("K134").Value = 3469.6362
$ws.Range("M134").Value = -934.6361999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1606.3334
$ws.Range("I31").Value = 1866.9231
$ws.Range("K31").Value = 1866.9231
$ws.Range("M31").Value = -1571.9231

$ws.Range("H34").Value = 1606.3334
$ws.Range("I34").Value = 1866.9231
$ws.Range("K34").Value = 1866.9231
$ws.Range("M34").Value = -1664.9231

$ws.Range("H58").Value = 1480.4615
$ws.Range("I58").Value = 1167
$ws.Range("J58").Value = 1982
$ws.Range("K58").Value = 1167
$ws.Range("L58").Value = 1982
$ws.Range("M58").Value = -964
$ws.Range("N58").Value = -2388

$ws.Range("H74").Value = 28666.666
$ws.Range("J74").Value = 33000
$ws.Range("L74").Value = 33000
$ws.Range("N74").Value = -34748

$ws.Range("H77").Value = 28666.666
$ws.Range("J77").Value = 33000
$ws.Range("L77").Value = 99000
$ws.Range("N77").Value = -107736

$ws.Range("H114").Value = 24866.4
$ws.Range("J114").Value = 24866.4
$ws.Range("L114").Value = 24866.4
$ws.Range("N114").Value = -33544.4

$ws.Range("H122").Value = 5439.3335
$ws.Range("I122").Value = 5886.1055
$ws.Range("K122").Value = 17658.3165
$ws.Range("M122").Value = -15208.3165

$ws.Range("H132").Value = 2573.3
$ws.Range("I132").Value = 2092.0625
$ws.Range("J132").Value = 4498.25
$ws.Range("K132").Value = 6276.1875
$ws.Range("L132").Value = 13494.75
$ws.Range("M132").Value = -3746.1875
$ws.Range("N132").Value = -18554.75

$ws.Range("H134").Value = 12501631
$ws.Range("I134").Value = 1701.8823
$ws.Range("K134").Value = 5105.6469
$ws.Range("M134").Value = -2570.6469

$ws.Range("H136").Value = 1480.4615
$ws.Range("I136").Value = 1167
$ws.Range("J136").Value = 1982
$ws.Range("K136").Value = 3501
$ws.Range("L136").Value = 5946
$ws.Range("M136").Value = -951
$ws.Range("N136").Value = -11046

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 91.125
$ws.Range("I2").Value = 38
$ws.Range("J2").Value = 123
$ws.Range("K2").Value = 228
$ws.Range("L2").Value = 738
$ws.Range("M2").Value = -115
$ws.Range("N2").Value = -964

$ws.Range("H109").Value = 67722.336
$ws.Range("I109").Value = 77556.53999999999
$ws.Range("K109").Value = 232669.62
$ws.Range("M109").Value = -231629.62

$ws.Range("H121").Value = 757.1539
$ws.Range("J121").Value = 1201.4286
$ws.Range("L121").Value = 3604.2858
$ws.Range("N121").Value = -6224.2858

$ws.Range("H122").Value = 1158
$ws.Range("J122").Value = 1178.5333
$ws.Range("L122").Value = 10606.7997
$ws.Range("N122").Value = -15506.7997

$ws.Range("H131").Value = 28573154
$ws.Range("J131").Value = 1909.2903
$ws.Range("L131").Value = 5727.8709
$ws.Range("N131").Value = -15807.8709

$ws.Range("H134").Value = 5000.758
$ws.Range("I134").Value = 2071.3333
$ws.Range("J134").Value = 7441.9443
$ws.Range("K134").Value = 6213.999899999999
$ws.Range("L134").Value = 22325.8329
$ws.Range("M134").Value = -1143.999899999999
$ws.Range("N134").Value = -32465.8329

$ws.Range("H140").Value = 24199.49
$ws.Range("I140").Value = 60848.53
$ws.Range("J140").Value = 3431.7
$ws.Range("K140").Value = 182545.59
$ws.Range("L140").Value = 10295.1
$ws.Range("M140").Value = -177365.59
$ws.Range("N140").Value = -20655.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 215.5
$ws.Range("I2").Value = 315.25
$ws.Range("J2").Value = 149
$ws.Range("K2").Value = 315.25
$ws.Range("L2").Value = 149
$ws.Range("M2").Value = -202.25
$ws.Range("N2").Value = -375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

$ws.Range("H132").Value = 2393.2903
$ws.Range("I132").Value = 1999.7
$ws.Range("J132").Value = 3108.9092
$ws.Range("K132").Value = 5999.1
$ws.Range("L132").Value = 9326.7276
$ws.Range("M132").Value = -3469.1
$ws.Range("N132").Value = -14386.7276

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 13890351
$ws.Range("I122").Value = 17858780
$ws.Range("J122").Value = 848.75
$ws.Range("K122").Value = 53576340
$ws.Range("L122").Value = 2546.25
$ws.Range("M122").Value = -53573890
$ws.Range("N122").Value = -7446.25

